# ICEBoard Bill of Materials - Rev 1.1 modifications
#
# The "Where Used" (designator) column had three stale/placeholder entries
# that get corrected to reflect the actual board reference designators:
#   - RGB LED row:       D1..D8  -> LED1..LED8
#   - Tactile Switch row: S5..S9  -> BLeft, BUp, BCenter, BDown, BRight
#   - Buzzer row:         B1      -> Buzzer
#
# Finally, the cursor/selection is left on F18 (just below the table),
# matching where the author's edit session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "LED1, LED2, LED3, LED4, LED5, LED6, LED7, LED8"
$ws.Range("F6").Value = "BLeft, BUp, BCenter, BDown, BRight"
$ws.Range("F15").Value = "Buzzer"

$ws.Range("F18").Select()
